# Refitting NCDEs to individual patients (for manuscript figure)
# Adds a new "Label" column (H) that flags each patient group as
# Control (0) or MDD (1), and refreshes the slightly-changed refit
# values for the "MDD 34" row under the 100-iteration block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header - copy the header formatting (bold, centered, bordered)
# from an existing header cell, then set the text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Label"

# Classification label per row: 0 = Control, 1 = MDD
$ws.Range("H2").Value  = 0   # Control 10 (100 iterations)
$ws.Range("H3").Value  = 0   # Control 5
$ws.Range("H4").Value  = 1   # MDD 38
$ws.Range("H5").Value  = 1   # MDD 9
$ws.Range("H6").Value  = 1   # MDD 49
$ws.Range("H7").Value  = 1   # MDD 26
$ws.Range("H8").Value  = 1   # MDD 34
$ws.Range("H9").Value  = 0   # Control 10 (200 iterations)
$ws.Range("H10").Value = 0   # Control 5
$ws.Range("H11").Value = 1   # MDD 38
$ws.Range("H12").Value = 1   # MDD 9
$ws.Range("H13").Value = 1   # MDD 49
$ws.Range("H14").Value = 1   # MDD 26
$ws.Range("H15").Value = 1   # MDD 34

# Refreshed refit values for MDD 34 under the 100-iteration block
$ws.Range("D8").Value = 0.3644065768291381
$ws.Range("E8").Value = 0.6355934231708619
